$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 35.160799
$ws.Range("H2").Value = 105.482397
$ws.Range("I2").Value = 0.2238945559395223
$ws.Range("J2").Value = 0.2238945559395223
$ws.Range("M2").Value = 0.01512566666666667
$ws.Range("N2").Value = 0.045377
$ws.Range("O2").Value = 0.01277774711566512
$ws.Range("P2").Value = 0.01277774711566512
$ws.Range("Q2").Value = 0.5318305254076666
$ws.Range("R2").Value = 4.786474728669
$ws.Range("S2").Value = 0.002860868016369354
$ws.Range("T2").Value = 0.002860868016369354

# Row 3
$ws.Range("G3").Value = 35.160799
$ws.Range("H3").Value = 105.482397
$ws.Range("I3").Value = 0.2238945559395223
$ws.Range("J3").Value = 0.2238945559395223
$ws.Range("M3").Value = 0.8534100000000001
$ws.Range("N3").Value = 2.56023
$ws.Range("O3").Value = 0.7209372919747741
$ws.Range("P3").Value = 0.7209372919747741
$ws.Range("Q3").Value = 30.00657747459
$ws.Range("R3").Value = 270.05919727131
$ws.Range("S3").Value = 0.1614139348469338
$ws.Range("T3").Value = 0.1614139348469338

# Row 4
$ws.Range("G4").Value = 35.160799
$ws.Range("H4").Value = 105.482397
$ws.Range("I4").Value = 0.2238945559395223
$ws.Range("J4").Value = 0.2238945559395223
$ws.Range("M4").Value = 0.1323446666666667
$ws.Range("N4").Value = 0.397034
$ws.Range("O4").Value = 0.1118011337973199
$ws.Range("P4").Value = 0.1118011337973199
$ws.Range("Q4").Value = 4.653344223388666
$ws.Range("R4").Value = 41.880098010498
$ws.Range("S4").Value = 0.02503166520508606
$ws.Range("T4").Value = 0.02503166520508606

# Row 5
$ws.Range("G5").Value = 35.160799
$ws.Range("H5").Value = 105.482397
$ws.Range("I5").Value = 0.2238945559395223
$ws.Range("J5").Value = 0.2238945559395223
$ws.Range("M5").Value = 0.1828703333333333
$ws.Range("N5").Value = 0.548611
$ws.Range("O5").Value = 0.154483827112241
$ws.Range("P5").Value = 0.154483827112241
$ws.Range("Q5").Value = 6.429867033396333
$ws.Range("R5").Value = 57.86880330056699
$ws.Range("S5").Value = 0.03458808787113313
$ws.Range("T5").Value = 0.03458808787113313

# Row 6
$ws.Range("I6").Value = 0.2050965007332699
$ws.Range("J6").Value = 0.2050965007332699
$ws.Range("M6").Value = 0.01512566666666667
$ws.Range("N6").Value = 0.045377
$ws.Range("O6").Value = 0.01277774711566512
$ws.Range("P6").Value = 0.01277774711566512
$ws.Range("Q6").Value = 0.4871783473543334
$ws.Range("R6").Value = 4.384605126189
$ws.Range("S6").Value = 0.002620671220677549
$ws.Range("T6").Value = 0.002620671220677549

# Row 7
$ws.Range("I7").Value = 0.2050965007332699
$ws.Range("J7").Value = 0.2050965007332699
$ws.Range("M7").Value = 0.8534100000000001
$ws.Range("N7").Value = 2.56023
$ws.Range("O7").Value = 0.7209372919747741
$ws.Range("P7").Value = 0.7209372919747741
$ws.Range("Q7").Value = 27.48724288179001
$ws.Range("R7").Value = 247.38518593611
$ws.Range("S7").Value = 0.1478617158321459
$ws.Range("T7").Value = 0.1478617158321459

# Row 8
$ws.Range("I8").Value = 0.2050965007332699
$ws.Range("J8").Value = 0.2050965007332699
$ws.Range("M8").Value = 0.1323446666666667
$ws.Range("N8").Value = 0.397034
$ws.Range("O8").Value = 0.1118011337973199
$ws.Range("P8").Value = 0.1118011337973199
$ws.Range("Q8").Value = 4.262652179815333
$ws.Range("R8").Value = 38.363869618338
$ws.Range("S8").Value = 0.02293002131984243
$ws.Range("T8").Value = 0.02293002131984243

# Row 9
$ws.Range("I9").Value = 0.2050965007332699
$ws.Range("J9").Value = 0.2050965007332699
$ws.Range("M9").Value = 0.1828703333333333
$ws.Range("N9").Value = 0.548611
$ws.Range("O9").Value = 0.154483827112241
$ws.Range("P9").Value = 0.154483827112241
$ws.Range("Q9").Value = 5.890019179769667
$ws.Range("R9").Value = 53.010172617927
$ws.Range("S9").Value = 0.03168409236060407
$ws.Range("T9").Value = 0.03168409236060407

# Row 10
$ws.Range("G10").Value = 51.53356533333334
$ws.Range("H10").Value = 154.600696
$ws.Range("I10").Value = 0.3281519491717758
$ws.Range("J10").Value = 0.3281519491717758
$ws.Range("M10").Value = 0.01512566666666667
$ws.Range("N10").Value = 0.045377
$ws.Range("O10").Value = 0.01277774711566512
$ws.Range("P10").Value = 0.01277774711566512
$ws.Range("Q10").Value = 0.779479531376889
$ws.Range("R10").Value = 7.015315782392
$ws.Range("S10").Value = 0.004193042622029547
$ws.Range("T10").Value = 0.004193042622029547

# Row 11
$ws.Range("G11").Value = 51.53356533333334
$ws.Range("H11").Value = 154.600696
$ws.Range("I11").Value = 0.3281519491717758
$ws.Range("J11").Value = 0.3281519491717758
$ws.Range("M11").Value = 0.8534100000000001
$ws.Range("N11").Value = 2.56023
$ws.Range("O11").Value = 0.7209372919747741
$ws.Range("P11").Value = 0.7209372919747741
$ws.Range("Q11").Value = 43.97925999112001
$ws.Range("R11").Value = 395.81333992008
$ws.Range("S11").Value = 0.2365769775921438
$ws.Range("T11").Value = 0.2365769775921438

# Row 12
$ws.Range("G12").Value = 51.53356533333334
$ws.Range("H12").Value = 154.600696
$ws.Range("I12").Value = 0.3281519491717758
$ws.Range("J12").Value = 0.3281519491717758
$ws.Range("M12").Value = 0.1323446666666667
$ws.Range("N12").Value = 0.397034
$ws.Range("O12").Value = 0.1118011337973199
$ws.Range("P12").Value = 0.1118011337973199
$ws.Range("Q12").Value = 6.820192526184889
$ws.Range("R12").Value = 61.381732735664
$ws.Range("S12").Value = 0.03668775997520504
$ws.Range("T12").Value = 0.03668775997520504

# Row 13
$ws.Range("G13").Value = 51.53356533333334
$ws.Range("H13").Value = 154.600696
$ws.Range("I13").Value = 0.3281519491717758
$ws.Range("J13").Value = 0.3281519491717758
$ws.Range("M13").Value = 0.1828703333333333
$ws.Range("N13").Value = 0.548611
$ws.Range("O13").Value = 0.154483827112241
$ws.Range("P13").Value = 0.154483827112241
$ws.Range("Q13").Value = 9.423960270361778
$ws.Range("R13").Value = 84.81564243325599
$ws.Range("S13").Value = 0.05069416898239751
$ws.Range("T13").Value = 0.05069416898239751

# Row 14
$ws.Range("G14").Value = 0.5955593333333333
$ws.Range("H14").Value = 1.786678
$ws.Range("I14").Value = 0.003792362411113143
$ws.Range("J14").Value = 0.003792362411113143
$ws.Range("M14").Value = 0.01512566666666667
$ws.Range("N14").Value = 0.045377
$ws.Range("O14").Value = 0.01277774711566512
$ws.Range("P14").Value = 0.01277774711566512
$ws.Range("Q14").Value = 0.009008231956222223
$ws.Range("R14").Value = 0.081074087606
$ws.Range("S14").Value = 0.0000484578478601578
$ws.Range("T14").Value = 0.0000484578478601578

# Row 15
$ws.Range("G15").Value = 0.5955593333333333
$ws.Range("H15").Value = 1.786678
$ws.Range("I15").Value = 0.003792362411113143
$ws.Range("J15").Value = 0.003792362411113143
$ws.Range("M15").Value = 0.8534100000000001
$ws.Range("N15").Value = 2.56023
$ws.Range("O15").Value = 0.7209372919747741
$ws.Range("P15").Value = 0.7209372919747741
$ws.Range("Q15").Value = 0.5082562906600001
$ws.Range("R15").Value = 4.57430661594
$ws.Range("S15").Value = 0.002734055486854834
$ws.Range("T15").Value = 0.002734055486854834

# Row 16
$ws.Range("G16").Value = 0.5955593333333333
$ws.Range("H16").Value = 1.786678
$ws.Range("I16").Value = 0.003792362411113143
$ws.Range("J16").Value = 0.003792362411113143
$ws.Range("M16").Value = 0.1323446666666667
$ws.Range("N16").Value = 0.397034
$ws.Range("O16").Value = 0.1118011337973199
$ws.Range("P16").Value = 0.1118011337973199
$ws.Range("Q16").Value = 0.07881910145022222
$ws.Range("R16").Value = 0.709371913052
$ws.Range("S16").Value = 0.0004239904173327873
$ws.Range("T16").Value = 0.0004239904173327873

# Row 17
$ws.Range("G17").Value = 0.5955593333333333
$ws.Range("H17").Value = 1.786678
$ws.Range("I17").Value = 0.003792362411113143
$ws.Range("J17").Value = 0.003792362411113143
$ws.Range("M17").Value = 0.1828703333333333
$ws.Range("N17").Value = 0.548611
$ws.Range("O17").Value = 0.154483827112241
$ws.Range("P17").Value = 0.154483827112241
$ws.Range("Q17").Value = 0.1089101338064444
$ws.Range("R17").Value = 0.9801912042579999
$ws.Range("S17").Value = 0.0005858586590653642
$ws.Range("T17").Value = 0.0005858586590653642

# Row 18
$ws.Range("G18").Value = 37.54313466666667
$ws.Range("H18").Value = 112.629404
$ws.Range("I18").Value = 0.2390646317443189
$ws.Range("J18").Value = 0.2390646317443189
$ws.Range("M18").Value = 0.01512566666666667
$ws.Range("N18").Value = 0.045377
$ws.Range("O18").Value = 0.01277774711566512
$ws.Range("P18").Value = 0.01277774711566512
$ws.Range("Q18").Value = 0.5678649405897778
$ws.Range("R18").Value = 5.110784465308
$ws.Range("S18").Value = 0.003054707408728517
$ws.Range("T18").Value = 0.003054707408728517

# Row 19
$ws.Range("G19").Value = 37.54313466666667
$ws.Range("H19").Value = 112.629404
$ws.Range("I19").Value = 0.2390646317443189
$ws.Range("J19").Value = 0.2390646317443189
$ws.Range("M19").Value = 0.8534100000000001
$ws.Range("N19").Value = 2.56023
$ws.Range("O19").Value = 0.7209372919747741
$ws.Range("P19").Value = 0.7209372919747741
$ws.Range("Q19").Value = 32.03968655588
$ws.Range("R19").Value = 288.35717900292
$ws.Range("S19").Value = 0.1723506082166959
$ws.Range("T19").Value = 0.1723506082166959

# Row 20
$ws.Range("G20").Value = 37.54313466666667
$ws.Range("H20").Value = 112.629404
$ws.Range("I20").Value = 0.2390646317443189
$ws.Range("J20").Value = 0.2390646317443189
$ws.Range("M20").Value = 0.1323446666666667
$ws.Range("N20").Value = 0.397034
$ws.Range("O20").Value = 0.1118011337973199
$ws.Range("P20").Value = 0.1118011337973199
$ws.Range("Q20").Value = 4.968633643081778
$ws.Range("R20").Value = 44.71770278773599
$ws.Range("S20").Value = 0.02672769687985362
$ws.Range("T20").Value = 0.02672769687985362

# Row 21
$ws.Range("G21").Value = 37.54313466666667
$ws.Range("H21").Value = 112.629404
$ws.Range("I21").Value = 0.2390646317443189
$ws.Range("J21").Value = 0.2390646317443189
$ws.Range("M21").Value = 0.1828703333333333
$ws.Range("N21").Value = 0.548611
$ws.Range("O21").Value = 0.154483827112241
$ws.Range("P21").Value = 0.154483827112241
$ws.Range("Q21").Value = 6.865525550871555
$ws.Range("R21").Value = 61.789729957844
$ws.Range("S21").Value = 0.03693161923904092
$ws.Range("T21").Value = 0.03693161923904092

